# Rectification feuille de route du 07 aout 2023
#
# - Row 5 / column C (afternoon slot for "CV de permanence" / "ANGLAIS (4h)")
#   gets a new note about the 14h check-in with the training manager.
# - The cell is styled like the other "highlight" cells in the sheet
#   (orange fill, kept border, centered + wrapped text) and the row is
#   made taller so the new text fits.
# - The active selection is moved to C8 (the cell the author ended up on).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New note in C5.
$c5 = $ws.Range("C5")
$c5.Value = "14h: bilan intermédiaire avec David (Manager de formation)"

# Match the orange highlight fill used elsewhere on the sheet (theme accent2 /
# RGB ED7D31), keep the existing border, center + wrap the text.
$c5.Interior.Color = 3243501
$c5.HorizontalAlignment = -4108  # xlCenter
$c5.VerticalAlignment = -4108    # xlCenter
$c5.WrapText = $true

# Row 5 needs to grow to fit the new wrapped text.
$ws.Rows("5").RowHeight = 45.75

# Move the selection to C8, matching the author's final cursor position.
$ws.Range("C8").Select() | Out-Null
